$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Scanner" to "Session"
$ws.Name = "Session"

# Delete row 33 (Student ID 888096 was a duplicate/erroneous manual log entry).
# This shifts all subsequent rows up by one, matching the new A1:F53 dimension.
$ws.Rows.Item(33).Delete()
